# Auto-generated Excel COM-interop edit script
# Applies the numeric LevePrice/LeveProfit recalculations described in the commit diff
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (one block per affected row).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 46000
$ws.Range("J3").Value = 46000
$ws.Range("L3").Value = 46000
$ws.Range("N3").Value = -46228

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 500
$ws.Range("I49").Value = 500
$ws.Range("K49").Value = 1500
$ws.Range("M49").Value = -1364

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2607.5
$ws.Range("J100").Value = 2475
$ws.Range("L100").Value = 2475
$ws.Range("N100").Value = -3557

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H102").Value = 46000
$ws.Range("J102").Value = 46000
$ws.Range("L102").Value = 46000
$ws.Range("N102").Value = -52490

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1177.4474
$ws.Range("I129").Value = 681.6
$ws.Range("J129").Value = 1354.5358
$ws.Range("K129").Value = 2044.8
$ws.Range("L129").Value = 4063.6074
$ws.Range("M129").Value = 2955.2
$ws.Range("N129").Value = -14063.6074

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2564.0781
$ws.Range("I138").Value = 3472.5833
$ws.Range("J138").Value = 2354.423
$ws.Range("K138").Value = 10417.7499
$ws.Range("L138").Value = 7063.268999999999
$ws.Range("M138").Value = -5277.749899999999
$ws.Range("N138").Value = -17343.269

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 364763.25
$ws.Range("I32").Value = 404697.12
$ws.Range("K32").Value = 404697.12
$ws.Range("M32").Value = -404410.12

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 11113611
$ws.Range("I61").Value = 41667996
$ws.Range("J61").Value = 2924.5454
$ws.Range("K61").Value = 41667996
$ws.Range("L61").Value = 2924.5454
$ws.Range("M61").Value = -41667784
$ws.Range("N61").Value = -3348.5454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 58333.332
$ws.Range("J109").Value = 58333.332
$ws.Range("L109").Value = 58333.332
$ws.Range("N109").Value = -61107.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1551.25
$ws.Range("I122").Value = 1102.75
$ws.Range("J122").Value = 1999.75
$ws.Range("K122").Value = 3308.25
$ws.Range("L122").Value = 5999.25
$ws.Range("M122").Value = -858.25
$ws.Range("N122").Value = -10899.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5082.8716
$ws.Range("I132").Value = 5740.2856
$ws.Range("J132").Value = 4315.8887
$ws.Range("K132").Value = 17220.8568
$ws.Range("L132").Value = 12947.6661
$ws.Range("M132").Value = -14690.8568
$ws.Range("N132").Value = -18007.6661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 11113611
$ws.Range("I136").Value = 41667996
$ws.Range("J136").Value = 2924.5454
$ws.Range("K136").Value = 125003988
$ws.Range("L136").Value = 8773.636200000001
$ws.Range("M136").Value = -125001438
$ws.Range("N136").Value = -13873.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H137").Value = 60780
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 53332
$ws.Range("J58").Value = 53332
$ws.Range("L58").Value = 53332
$ws.Range("N58").Value = -53920

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H59").Value = 98780
$ws.Range("J59").Value = 98780
$ws.Range("L59").Value = 98780
$ws.Range("N59").Value = -100474

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4489.875
$ws.Range("I134").Value = 5228
$ws.Range("J134").Value = 3751.75
$ws.Range("K134").Value = 15684
$ws.Range("L134").Value = 11255.25
$ws.Range("M134").Value = -13149
$ws.Range("N134").Value = -16325.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3603.4905
$ws.Range("I31").Value = 1079.5333
$ws.Range("J31").Value = 6895.609
$ws.Range("K31").Value = 1079.5333
$ws.Range("L31").Value = 6895.609
$ws.Range("M31").Value = -784.5333000000001
$ws.Range("N31").Value = -7485.609

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3603.4905
$ws.Range("I34").Value = 1079.5333
$ws.Range("J34").Value = 6895.609
$ws.Range("K34").Value = 1079.5333
$ws.Range("L34").Value = 6895.609
$ws.Range("M34").Value = -877.5333000000001
$ws.Range("N34").Value = -7299.609

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 99090.336
$ws.Range("J63").Value = 99090.336
$ws.Range("L63").Value = 99090.336
$ws.Range("N63").Value = -100462.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H66").Value = 99090.336
$ws.Range("J66").Value = 99090.336
$ws.Range("L66").Value = 297271.008
$ws.Range("N66").Value = -304135.008

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 28166
$ws.Range("J68").Value = 28166
$ws.Range("L68").Value = 28166
$ws.Range("N68").Value = -29664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 28166
$ws.Range("J71").Value = 28166
$ws.Range("L71").Value = 84498
$ws.Range("N71").Value = -91986

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1894.6842
$ws.Range("I99").Value = 1000
$ws.Range("K99").Value = 1000
$ws.Range("M99").Value = 498

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1894.6842
$ws.Range("I126").Value = 1000
$ws.Range("K126").Value = 3000
$ws.Range("M126").Value = -530

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 10419002
$ws.Range("I132").Value = 1072.4
$ws.Range("J132").Value = 15154424
$ws.Range("K132").Value = 3217.2
$ws.Range("L132").Value = 45463272
$ws.Range("M132").Value = -687.2000000000003
$ws.Range("N132").Value = -45468332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 23809900
$ws.Range("I34").Value = 166
$ws.Range("J34").Value = 33333792
$ws.Range("K34").Value = 498
$ws.Range("L34").Value = 100001376
$ws.Range("M34").Value = -414
$ws.Range("N34").Value = -100001544

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3655.5557
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 3655.5557
$ws.Range("K39").Value = 0
$ws.Range("L39").ClearContents()
$ws.Range("M39").Value = 10966.6671
$ws.Range("N39").Value = -11554.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1900
$ws.Range("J55").Value = 1900
$ws.Range("L55").Value = 5700
$ws.Range("N55").Value = -6054

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2461.8489
$ws.Range("I68").Value = 3354.743
$ws.Range("J68").Value = 1849.0785
$ws.Range("K68").Value = 10064.229
$ws.Range("L68").Value = 5547.235500000001
$ws.Range("M68").Value = -9253.228999999999
$ws.Range("N68").Value = -7169.235500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2461.8489
$ws.Range("I71").Value = 3354.743
$ws.Range("J71").Value = 1849.0785
$ws.Range("K71").Value = 30192.687
$ws.Range("L71").Value = 16641.7065
$ws.Range("M71").Value = -26136.687
$ws.Range("N71").Value = -24753.7065

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2138.151
$ws.Range("I107").Value = 375.53845
$ws.Range("J107").Value = 2711
$ws.Range("K107").Value = 1126.61535
$ws.Range("L107").Value = 8133
$ws.Range("M107").Value = 793.38465
$ws.Range("N107").Value = -11973

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1009.82855
$ws.Range("J113").Value = 1469.7059
$ws.Range("L113").Value = 4409.1177
$ws.Range("N113").Value = -8749.117699999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1426.0435
$ws.Range("I129").Value = 684.875
$ws.Range("J129").Value = 1821.3334
$ws.Range("K129").Value = 2054.625
$ws.Range("L129").Value = 5464.0002
$ws.Range("M129").Value = 2945.375
$ws.Range("N129").Value = -15464.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 863.3103599999999
$ws.Range("I131").Value = 342.72726
$ws.Range("J131").Value = 1181.4445
$ws.Range("K131").Value = 1028.18178
$ws.Range("L131").Value = 3544.3335
$ws.Range("M131").Value = 4011.81822
$ws.Range("N131").Value = -13624.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 6231.294
$ws.Range("I137").Value = 8716.125
$ws.Range("J137").Value = 4022.5557
$ws.Range("K137").Value = 26148.375
$ws.Range("L137").Value = 12067.6671
$ws.Range("M137").Value = -21048.375
$ws.Range("N137").Value = -22267.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 1894.3572
$ws.Range("I138").Value = 1308.6364
$ws.Range("J138").Value = 4042
$ws.Range("K138").Value = 3925.9092
$ws.Range("L138").Value = 12126
$ws.Range("M138").Value = 1214.0908
$ws.Range("N138").Value = -22406

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2360.8333
$ws.Range("I140").Value = 1332.2727
$ws.Range("J140").Value = 3977.1428
$ws.Range("K140").Value = 3996.8181
$ws.Range("L140").Value = 11931.4284
$ws.Range("M140").Value = 1183.1819
$ws.Range("N140").Value = -22291.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3728.75
$ws.Range("I132").Value = 3532.5
$ws.Range("J132").Value = 3990.4167
$ws.Range("K132").Value = 10597.5
$ws.Range("L132").Value = 11971.2501
$ws.Range("M132").Value = -8067.5
$ws.Range("N132").Value = -17031.2501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5611244
$ws.Range("I132").Value = 2494.8064
$ws.Range("J132").Value = 13890827
$ws.Range("K132").Value = 7484.4192
$ws.Range("L132").Value = 41672481
$ws.Range("M132").Value = -4954.4192
$ws.Range("N132").Value = -41677541

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("N133").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3081.2559
$ws.Range("I136").Value = 2809.6775
$ws.Range("K136").Value = 8429.032499999999
$ws.Range("M136").Value = -5879.032499999999
